$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.080.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5298'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.34%  '

$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2617'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06336'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07765'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.494'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.656.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5482'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅8144'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.20'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.106.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.57%  '

$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.548'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.022'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.20%  '

$ws.Range("E23").Value = '  -0.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '140.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1243'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.271'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.47%  '

$ws.Range("E28").Value = '  +1.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05945'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.277'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.508'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.235'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.559'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9486'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.411'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.768'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5642'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01610'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.838'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8456'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.011.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.801.93'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.61%  '

$ws.Range("E47").Value = '  -8.54%  '

$ws.Range("E48").Value = '  +1.31%  '

$ws.Range("E49").Value = '  -0.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.468'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.716'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.42%  '
